$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 3) with the new item
$ws.Range("A3").Value = "Sauce Labs Bike Light"
$ws.Range("B3").Value = "A red light isn't the desired state in testing but it sure helps when riding your bike at night. Water-resistant with 3 lighting modes, 1 AAA battery included."
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '$9.99'

# Update the selection to match the target state
$ws.Range("E7").Select()
